# Re-save the workbook (as Excel does when it opens a foreign-generated
# .xlsx and writes it back out), selecting the used range on Sheet1 so the
# saved sheetView records activeCell="A4" / sqref="A4:H30" - matching what
# happened when the terminal data was reviewed/tested in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
[void]$ws.Range("A4:H30").Select()
